$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A31").Value = $ws.Range("A30").Value
$ws.Range("B31").Value = "19. Remove Nth Node From End of List"
$ws.Range("C31").Value = 'Use 2 ptrs L & R, initialize both = head, use while loop to move R by "n-1" places to the right of L' + "`n" + 'Initalize a prev = null' + "`n" + 'Now move both using while(R.next!=null) since both will maintain a fixed distance, after the end of the loop L will point to the node to be deleted. The "prev" pointer points to the node before L. Using prev delete L'

$ws.Range("A31").Style = $ws.Range("A30").Style
$ws.Range("B31").Style = $ws.Range("B17").Style
$ws.Range("C31").Style = $ws.Range("C30").Style

$ws.Rows.Item(31).RowHeight = $ws.Rows.Item(17).RowHeight

$ws.Range("C35").Select()
